$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking values stored as literal text.
# Prefixing the new value with a leading apostrophe tells Excel to keep it
# as text (quotePrefix) instead of auto-converting it to a Number, matching
# how the source data is stored (inline/shared text strings).

$ws.Range("D2").Value = "'243.11"
$ws.Range("D3").Value = "'22.97"
$ws.Range("D4").Value = "'5.381"
$ws.Range("D5").Value = "'0.05907"
$ws.Range("D6").Value = "'3.456"
$ws.Range("D7").Value = "'6.550"
$ws.Range("D8").Value = "'0.8107"
$ws.Range("D9").Value = "'0.9156"
$ws.Range("D10").Value = "'0.1413"
$ws.Range("D11").Value = "'0.07422"
$ws.Range("D12").Value = "'0.03284"
$ws.Range("D13").Value = "'0.03069"
$ws.Range("D14").Value = "'0.09345"
$ws.Range("D15").Value = "'3.855"
$ws.Range("D16").Value = "'0.001571"
$ws.Range("D17").Value = "'0.04673"
$ws.Range("D18").Value = "'0.0005921"
$ws.Range("D19").Value = "'0.006008"
$ws.Range("D20").Value = "'0.001296"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("D21").Value = "'0.004927"
$ws.Range("D22").Value = "'0.00009501"
$ws.Range("D24").Value = "'2.151"
$ws.Range("D25").Value = "'0.3242"
$ws.Range("D27").Value = "'0.0002900"
$ws.Range("D40").Value = "'0.03951"
$ws.Range("D41").Value = "'0.006192"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D43").Value = "'0.002617"
$ws.Range("D44").Value = "'0.008085"
$ws.Range("D45").Value = "'0.00005187"
$ws.Range("D47").Value = "'0.7901"
$ws.Range("D48").Value = "'0.002283"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
